# The author inserted one new weekly price record for "Albahaca" (Vega
# Modelo de Temuco) ahead of the existing row 316, which pushes every
# following record down by one row (old row 316 -> new row 317, ...,
# old row 425 -> new row 426). The sheet's used range grows from
# A1:R425 to A1:R426 as a result.
#
# Replicate that with a real row insert (so every row below shifts down
# and formatting/styles on column D - the date column - carry over the
# way Excel does it), then populate the newly freed row 316 with the
# new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 316; rows 316:425 shift down to 317:426.
$ws.Rows("316:316").Insert()

# Fill in the new record at row 316.
$ws.Cells.Item(316, 1).Value = 10
$ws.Cells.Item(316, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(316, 3).Value = "La Araucanía"
$ws.Cells.Item(316, 4).Value = 45146
$ws.Cells.Item(316, 5).Value = 9
$ws.Cells.Item(316, 6).Value = 100112052
$ws.Cells.Item(316, 7).Value = "Albahaca"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 35
$ws.Cells.Item(316, 11).Value = 6000
$ws.Cells.Item(316, 12).Value = 6000
$ws.Cells.Item(316, 13).Value = 6000
$ws.Cells.Item(316, 14).Value = "`$/paquete"
$ws.Cells.Item(316, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(316, 16).Value = 6000
$ws.Cells.Item(316, 17).Value = 1
$ws.Cells.Item(316, 18).Value = "Hortaliza"
